$d = $word.ActiveDocument

function Find-Paragraph($doc, [string]$startsWith) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.StartsWith($startsWith)) {
            return $p
        }
    }
    return $null
}

function Apply-Breaks($para, $pairs) {
    foreach ($pair in $pairs) {
        $search = $pair[0]
        $replace = $pair[1]
        $r = $para.Range
        $found = $r.Find.Execute($search, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
        if (-not $found) {
            Write-Host "WARNING: not found ->" $search
        }
    }
}

$pairsPt = @(
    @("cos de ecologia da paisagem.", "cos de ecologia da paisagem.^l"),
    @(" reabilitação e recuperação;", " reabilitação e recuperação;^l"),
    @("emas aquáticos e terrestres;", "emas aquáticos e terrestres;^l"),
    @(" Restauração vs. Recuperação", " Restauração vs. Recuperação^l"),
    @("eiras e regeneração natural;", "eiras e regeneração natural;^l"),
    @("o natural. Saída para Campo;", "o natural. Saída para Campo;^l"),
    @("referentes à saída de campo;", "referentes à saída de campo;^l"),
    @("es radiculares específicas);", "es radiculares específicas);^l"),
    @("os de restauração ecológica;", "os de restauração ecológica;^l"),
    @("ção às Mudanças Climáticas);", "ção às Mudanças Climáticas);^l"),
    @("to de restauração ecológica.", "to de restauração ecológica.^l"),
)

$pairsEn = @(
    @("entals of landscape ecology.", "entals of landscape ecology.^l"),
    @("rehabilitation and recovery;", "rehabilitation and recovery;^l"),
    @("tic and terrestrial systems;", "tic and terrestrial systems;^l"),
    @("- Restoration vs. Recovery", "- Restoration vs. Recovery^l"),
    @("cs and natural regeneration;", "cs and natural regeneration;^l"),
    @("ion. Departure to the Field;", "ion. Departure to the Field;^l"),
    @("ata relating to field trips;", "ata relating to field trips;^l"),
    @(" (specific root properties);", " (specific root properties);^l"),
    @("ogical restoration projects;", "ogical restoration projects;^l"),
    @("te Change Adaptation Plans);", "te Change Adaptation Plans);^l"),
    @("logical restoration project.", "logical restoration project.^l"),
)

$pairsBiblio = @(
    @("Bibliografia básica:", "Bibliografia básica:^l^l"),
    @("cina de Textos, 432 p. 2015.", "cina de Textos, 432 p. 2015.^l^l"),
    @(". Island Press, 336 p. 2013.", ". Island Press, 336 p. 2013.^l^l"),
    @("ina, PR: Planta, 300 p. 2013", "ina, PR: Planta, 300 p. 2013^l^l"),
    @("aturais. FEPAF, 340 p. 2003.", "aturais. FEPAF, 340 p. 2003.^l^l"),
    @("s degradados. 2ª Ed. Viçosa,", "s degradados. 2ª Ed. Viçosa,^l"),
    @(" Viçosa (UFV), 376 p., 2015.", " Viçosa (UFV), 376 p., 2015.^l^l"),
    @("l Restoration, 584 p., 2013.", "l Restoration, 584 p., 2013.^l^l^l"),
    @("Bibliografia complementar:", "Bibliografia complementar:^l^l"),
    @("stems. Forests, 14(7), 1442.", "stems. Forests, 14(7), 1442.^l^l"),
    @("ditora Forense, 176 p. 2013.", "ditora Forense, 176 p. 2013.^l^l"),
    @("gia & Ambiente, 5(1), 40-48.", "gia & Ambiente, 5(1), 40-48.^l^l"),
)

$paraPt = Find-Paragraph $d "- Estrutura de populações"
if ($paraPt -eq $null) { Write-Host "ERROR: Portuguese Programa paragraph not found" } else { Apply-Breaks $paraPt $pairsPt }

$paraEn = Find-Paragraph $d "- Structure of populations"
if ($paraEn -eq $null) { Write-Host "ERROR: English Programa paragraph not found" } else { Apply-Breaks $paraEn $pairsEn }

$paraBiblio = Find-Paragraph $d "Bibliografia básica:"
if ($paraBiblio -eq $null) { Write-Host "ERROR: Bibliografia paragraph not found" } else { Apply-Breaks $paraBiblio $pairsBiblio }

Write-Host "Done."

